$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = '展览'; Cell = 'F2'; Value = 143 }
    @{ Sheet = '展览'; Cell = 'F3'; Value = 963 }
    @{ Sheet = '展览'; Cell = 'F5'; Value = 2940 }
    @{ Sheet = '展览'; Cell = 'F6'; Value = 789 }
    @{ Sheet = '展览'; Cell = 'F7'; Value = 594 }
    @{ Sheet = '展览'; Cell = 'F9'; Value = 418 }
    @{ Sheet = '展览'; Cell = 'F10'; Value = 657 }
    @{ Sheet = '展览'; Cell = 'F12'; Value = 516 }
    @{ Sheet = '展览'; Cell = 'F14'; Value = 2156 }
    @{ Sheet = '展览'; Cell = 'F15'; Value = 1260 }
    @{ Sheet = '展览'; Cell = 'F16'; Value = 744 }
    @{ Sheet = '展览'; Cell = 'F17'; Value = 17 }
    @{ Sheet = '展览'; Cell = 'F19'; Value = 2666 }
    @{ Sheet = '展览'; Cell = 'E23'; Value = '2024.04.20 10:00-04.21 17:00' }
    @{ Sheet = '展览'; Cell = 'F23'; Value = 529 }
    @{ Sheet = '展览'; Cell = 'I23'; Value = '//i2.hdslb.com/bfs/openplatform/202404/I2fHOXWb1712026178009.jpeg' }
    @{ Sheet = '展览'; Cell = 'F25'; Value = 605 }
    @{ Sheet = '展览'; Cell = 'F27'; Value = 21 }
    @{ Sheet = '展览'; Cell = 'F32'; Value = 117 }
    @{ Sheet = '展览'; Cell = 'F34'; Value = 4670 }
    @{ Sheet = '展览'; Cell = 'F35'; Value = 239 }
    @{ Sheet = '展览'; Cell = 'F36'; Value = 22 }
    @{ Sheet = '演出'; Cell = 'F23'; Value = 274 }
    @{ Sheet = '演出'; Cell = 'F25'; Value = 299 }
    @{ Sheet = '演出'; Cell = 'F27'; Value = 136 }
    @{ Sheet = '演出'; Cell = 'F31'; Value = 22 }
    @{ Sheet = '演出'; Cell = 'F36'; Value = 539 }
    @{ Sheet = '本地生活'; Cell = 'F6'; Value = 243 }
    @{ Sheet = '本地生活'; Cell = 'F7'; Value = 251 }
    @{ Sheet = '全部类型'; Cell = 'F5'; Value = 143 }
    @{ Sheet = '全部类型'; Cell = 'F6'; Value = 243 }
    @{ Sheet = '全部类型'; Cell = 'F7'; Value = 963 }
    @{ Sheet = '全部类型'; Cell = 'F9'; Value = 2940 }
    @{ Sheet = '全部类型'; Cell = 'F10'; Value = 789 }
    @{ Sheet = '全部类型'; Cell = 'F11'; Value = 594 }
    @{ Sheet = '全部类型'; Cell = 'F13'; Value = 418 }
    @{ Sheet = '全部类型'; Cell = 'F14'; Value = 657 }
    @{ Sheet = '全部类型'; Cell = 'F16'; Value = 516 }
    @{ Sheet = '全部类型'; Cell = 'F21'; Value = 2156 }
    @{ Sheet = '全部类型'; Cell = 'F22'; Value = 1260 }
    @{ Sheet = '全部类型'; Cell = 'F23'; Value = 744 }
    @{ Sheet = '全部类型'; Cell = 'F27'; Value = 2666 }
    @{ Sheet = '全部类型'; Cell = 'E32'; Value = '2024.04.20 10:00-04.21 17:00' }
    @{ Sheet = '全部类型'; Cell = 'F32'; Value = 529 }
    @{ Sheet = '全部类型'; Cell = 'I32'; Value = '//i2.hdslb.com/bfs/openplatform/202404/I2fHOXWb1712026178009.jpeg' }
    @{ Sheet = '全部类型'; Cell = 'F35'; Value = 251 }
    @{ Sheet = '全部类型'; Cell = 'F37'; Value = 605 }
    @{ Sheet = '全部类型'; Cell = 'F38'; Value = 605 }
    @{ Sheet = '全部类型'; Cell = 'F39'; Value = 274 }
    @{ Sheet = '全部类型'; Cell = 'F42'; Value = 299 }
    @{ Sheet = '全部类型'; Cell = 'F47'; Value = 4670 }
    @{ Sheet = '全部类型'; Cell = 'F48'; Value = 239 }
    @{ Sheet = '全部类型'; Cell = 'F50'; Value = 539 }
)

foreach ($change in $changes) {
    $ws = $wb.Worksheets.Item($change.Sheet)
    $ws.Range($change.Cell).Value = $change.Value
}
